$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column A (ID) - shifts everything right by one
$ws.Columns("A").Insert()

# Insert new column K (Description) - shifts Easy..Sunkid Moving Carpet right by one
$ws.Columns("K").Insert()

# --- Headers (row 1) ---
# New column A has no header style yet (it was a brand-new empty column);
# copy the bold/bordered header formatting from B1 so A1 matches the rest
# of the header row (K1 already inherited it from its neighbours).
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = "ID"
$ws.Range("K1").Value = "Description"

# --- Row 2: KitzSki ---
$ws.Range("A2").Value = "kitzski-kitzbuehelkirchberg"
$ws.Range("B2").Value = "Ski resort KitzSki – Kitzbühel/​Kirchberg"
$ws.Range("K2").Value = "Kitzbühel - the legend: where skiing was ‘invented’, a venue for regular ski races since 1895 and home to the Streif, for decades the arena for the most famous international race in the Winter World Cup. But it is not just Kitzbühel's fantastic selection of slopes that make it such an irresistible destination for many visitors. The resort's unique blend of traditional and contemporary attracts skiers and celebrities from all over the world. An above-average number of luxury hotels, the highest concentration of award-winning restaurants in Tyrol and an amazing evening scene including numerous events make the winter sports resort and town of Kitzbühel Austria's meeting place for snow society throughout the entire winter."

# --- Row 3: Zillertal Arena ---
$ws.Range("A3").Value = "zillertal-arena-zell-am-zillergerloskoenigsleitenhochkrimml"
$ws.Range("B3").Value = "Ski resort Zillertal Arena – Zell am Ziller/​Gerlos/​Königsleiten/​Hochkrimml"
$ws.Range("K3").Value = "In the Zillertal Arena ski resort, extending from the Zillertal valley in Tyrol to Pinzgau in the Province of Salzburg, winter sports enthusiasts can enjoy an unbeatable variety of ski slopes. Varied, family-friendly and snow reliable - the Zillertal Arena offers recreation and fun for all ages. The winter sports resort's facilities range from SkiMovie routes and photo points to free WiFi hotspots throughout the entire ski resort. You can round off an action-packed day on the slopes with a ride on the Arena Coaster (alpine roller coaster) or a trip down the 7-kilometre Gerlosstein toboggan run."

# --- Row 4: San Martino di Castrozza ---
$ws.Range("A4").Value = "san-martino-di-castrozza"
$ws.Range("B4").Value = "Ski resort San Martino di Castrozza"
$ws.Range("K4").Value = "The ski resort of San Martino di Castrozza is situated in a beautiful Dolomite landscape. The amazing views of the Pale di San Martino (Pala group) and the traditional hospitality make a winter holiday in the ski resort of San Martino di Castrozza something special. The well-prepared slopes, most of which have snow-making capabilities, cater to all difficulty levels, with different gradients and offer lots of variety for fun skiing and snowboarding. There is a very good selection for both beginners and experts. There are two separate areas for children with moving carpets and practice lifts. Beautiful cross-country trails complete the varied offering in the winter sports resort of San Martino di Castrozza."

# --- Row 5: Paganella - Andalo ---
$ws.Range("A5").Value = "paganella-andalo"
$ws.Range("B5").Value = "Ski resort Paganella – Andalo"
$ws.Range("K5").Value = "It is possible to enter the ski area on the Paganella directly from Andalo and Fai della Paganella. Lifts provide access to slopes of all difficulty levels. Some of the slopes for skiing and snowboarding are above the treeline and others are below. Different children's areas and beginners’ areas are distributed throughout the ski resort. The winter sports resort of Andalo also offers a wide range of other activities including cross-country skiing, snow-shoeing and winter hiking, tobogganing and much more. The view of the Dolomites and Lake Garda is magnificent."

